$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 716, pushing the existing rows 716-813 down to 718-815.
$ws.Rows.Item(716).Insert()
$ws.Rows.Item(716).Insert()

# New row 716
$ws.Range("A716").Value = 3
$ws.Range("B716").Value = "Femacal de La Calera"
$ws.Range("C716").Value = "Coquimbo"
$ws.Range("D716").Value = 45131
$ws.Range("E716").Value = 5
$ws.Range("F716").Value = 100112003
$ws.Range("G716").Value = "Ajo"
$ws.Range("H716").Value = "Chino"
$ws.Range("I716").Value = "Primera"
$ws.Range("J716").Value = 45
$ws.Range("K716").Value = 18000
$ws.Range("L716").Value = 18000
$ws.Range("M716").Value = 18000
$ws.Range("N716").Value = "$/caja 10 kilos"
$ws.Range("O716").Value = "China"
$ws.Range("P716").Value = 1800
$ws.Range("Q716").Value = 10
$ws.Range("R716").Value = "Hortaliza"

# New row 717
$ws.Range("A717").Value = 3
$ws.Range("B717").Value = "Femacal de La Calera"
$ws.Range("C717").Value = "Coquimbo"
$ws.Range("D717").Value = 45131
$ws.Range("E717").Value = 5
$ws.Range("F717").Value = 100112003
$ws.Range("G717").Value = "Ajo"
$ws.Range("H717").Value = "Chino"
$ws.Range("I717").Value = "Primera"
$ws.Range("J717").Value = 78
$ws.Range("K717").Value = 22000
$ws.Range("L717").Value = 22500
$ws.Range("M717").Value = 22256
$ws.Range("N717").Value = "$/malla 10 kilos"
$ws.Range("O717").Value = "China"
$ws.Range("P717").Value = 2226
$ws.Range("Q717").Value = 10
$ws.Range("R717").Value = "Hortaliza"
